## Update 2p3. Added templates for formula student suspension, torque
## vectoring, four-wheel steering.
##
## Concretely: add a new "FSAE_Achilles" sheet (cloned from the
## "Bus_Makhulu" template) at the end of the workbook, populate it with
## the new vehicle's body data, and nudge a couple of leftover cursor
## positions on the other sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet by duplicating the "Bus_Makhulu" template and
#    moving the duplicate to the end of the tab strip.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("Bus_Makhulu")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Reflection.Missing]::Value, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "FSAE_Achilles"

# ---------------------------------------------------------------------
# 2. Fill in the new vehicle's numbers on the new sheet only.
# ---------------------------------------------------------------------
$ws.Range("H3").Value2 = "FSAE_Achilles"

$ws.Range("F6").Value2 = -1.53

$ws.Range("F7").Value2 = -0.8
$ws.Range("G7").Value2 = 0
$ws.Range("H7").Value2 = 0.289

$ws.Range("F8").Value2 = -1

$ws.Range("F9").Value2 = 0.25
$ws.Range("H9").Value2 = 0.403

$ws.Range("F10").Value2 = -1.75
$ws.Range("H10").Value2 = 0.403

$ws.Range("H11").Formula = "=0.619*2+0.2"

$ws.Range("H12").Value2 = 165

$ws.Range("F13").Value2 = 43
$ws.Range("G13").Value2 = 192
$ws.Range("H13").Value2 = 206

# Style touch-ups that differ from the Bus_Makhulu template: row 7's
# format spreads into G7/H11/H12, and row13 picks up the row7 look.
$ws.Range("G11").Copy()
$ws.Range("H12").PasteSpecial(-4122)

$ws.Range("F7").Copy()
$ws.Range("F13:H13").PasteSpecial(-4122)

# Bus_Makhulu's leftover "guesses" helper note in K12 doesn't apply here.
$ws.Range("K12").ClearContents()
$ws.Range("K12").ClearFormats()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Misc cursor-position housekeeping left behind on the other sheets.
# ---------------------------------------------------------------------
$sedanHamba = $wb.Worksheets.Item("Sedan_Hamba")
$sedanHamba.Activate()
$sedanHamba.Range("E22").Select()

$sedanHambaLG = $wb.Worksheets.Item("Sedan_HambaLG")
$sedanHambaLG.Activate()
$sedanHambaLG.Range("H12").Select()

# ---------------------------------------------------------------------
# 4. Leave the new sheet active/selected, as it is the one being worked
#    on.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("G27").Select()
